$wb = $excel.ActiveWorkbook

# --- ValidLoginsheet: update the stored email and move the selection there ---
$wsValid = $wb.Worksheets.Item("ValidLoginsheet")
$wsValid.Range("A2").Value = "adityapawar123@yopmail.com"
$wsValid.Range("A2").Select()

# --- signupdata: make it the active sheet, with C2 selected ---
$wsSignup = $wb.Worksheets.Item("signupdata")
$wsSignup.Activate()
$wsSignup.Range("C2").Select()
